# Auto-generated edit script: updates market/profit data cells
# across multiple crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -3082
$ws.Range("H112").Value = 25002456
$ws.Range("J112").Value = 26318324
$ws.Range("L112").Value = 78954972
$ws.Range("N112").Value = -78957188
$ws.Range("H113").Value = 1500
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -8008
$ws.Range("H132").Value = 2987.087
$ws.Range("I132").Value = 3105.0244
$ws.Range("K132").Value = 9315.073199999999
$ws.Range("M132").Value = -6785.073199999999
$ws.Range("H137").Value = 14286667
$ws.Range("I137").Value = 823.62744
$ws.Range("J137").Value = 52632876
$ws.Range("K137").Value = 2470.88232
$ws.Range("L137").Value = 157898628
$ws.Range("M137").Value = 79.11768000000029
$ws.Range("N137").Value = -157903728

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2470339.2
$ws.Range("I61").Value = 3473474
$ws.Range("J61").Value = 1084.3077
$ws.Range("K61").Value = 3473474
$ws.Range("L61").Value = 1084.3077
$ws.Range("M61").Value = -3473262
$ws.Range("N61").Value = -1508.3077
$ws.Range("H74").Value = 14290815
$ws.Range("I74").Value = 19231720
$ws.Range("J74").Value = 17086.334
$ws.Range("K74").Value = 19231720
$ws.Range("L74").Value = 17086.334
$ws.Range("M74").Value = -19230846
$ws.Range("N74").Value = -18834.334
$ws.Range("H77").Value = 14290815
$ws.Range("I77").Value = 19231720
$ws.Range("J77").Value = 17086.334
$ws.Range("K77").Value = 96158600
$ws.Range("L77").Value = 85431.67
$ws.Range("M77").Value = -96154232
$ws.Range("N77").Value = -94167.67
$ws.Range("I132").Value = 11224908
$ws.Range("J132").Value = 93097.91
$ws.Range("K132").Value = 33674724
$ws.Range("L132").Value = 279293.73
$ws.Range("M132").Value = -33672194
$ws.Range("N132").Value = -284353.73
$ws.Range("H136").Value = 2470339.2
$ws.Range("I136").Value = 3473474
$ws.Range("J136").Value = 1084.3077
$ws.Range("K136").Value = 10420422
$ws.Range("L136").Value = 3252.9231
$ws.Range("M136").Value = -10417872
$ws.Range("N136").Value = -8352.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 641561.3
$ws.Range("I107").Value = 829785.1
$ws.Range("J107").Value = 1600.4
$ws.Range("K107").Value = 829785.1
$ws.Range("L107").Value = 1600.4
$ws.Range("M107").Value = -827865.1
$ws.Range("N107").Value = -5440.4
$ws.Range("H134").Value = 3591864
$ws.Range("I134").Value = 4190318.5
$ws.Range("J134").Value = 1137.375
$ws.Range("K134").Value = 12570955.5
$ws.Range("L134").Value = 3412.125
$ws.Range("M134").Value = -12568420.5
$ws.Range("N134").Value = -8482.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4394602.5
$ws.Range("I31").Value = 1164.4
$ws.Range("J31").Value = 26361794
$ws.Range("K31").Value = 1164.4
$ws.Range("L31").Value = 26361794
$ws.Range("M31").Value = -869.4000000000001
$ws.Range("N31").Value = -26362384
$ws.Range("H34").Value = 4394602.5
$ws.Range("I34").Value = 1164.4
$ws.Range("J34").Value = 26361794
$ws.Range("K34").Value = 1164.4
$ws.Range("L34").Value = 26361794
$ws.Range("M34").Value = -962.4000000000001
$ws.Range("N34").Value = -26362198
$ws.Range("H58").Value = 1321.3383
$ws.Range("I58").Value = 979.2593000000001
$ws.Range("K58").Value = 979.2593000000001
$ws.Range("M58").Value = -776.2593000000001
$ws.Range("H99").Value = 1484.3334
$ws.Range("I99").Value = 1528.3636
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 1528.3636
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = -30.36359999999991
$ws.Range("N99").Value = -3996
$ws.Range("H126").Value = 1484.3334
$ws.Range("I126").Value = 1528.3636
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 4585.0908
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -2115.0908
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 1740.6428
$ws.Range("I132").Value = 1316.05
$ws.Range("K132").Value = 3948.15
$ws.Range("M132").Value = -1418.15
$ws.Range("H134").Value = 1383.8
$ws.Range("I134").Value = 1590.2
$ws.Range("J134").Value = 867.8
$ws.Range("K134").Value = 4770.6
$ws.Range("L134").Value = 2603.4
$ws.Range("M134").Value = -2235.6
$ws.Range("N134").Value = -7673.4
$ws.Range("H136").Value = 1321.3383
$ws.Range("I136").Value = 979.2593000000001
$ws.Range("K136").Value = 2937.7779
$ws.Range("M136").Value = -387.7779

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 2407.2727
$ws.Range("I123").Value = 1245
$ws.Range("J123").Value = 3071.4285
$ws.Range("K123").Value = 3735
$ws.Range("L123").Value = 9214.2855
$ws.Range("M123").Value = -1285
$ws.Range("N123").Value = -14114.2855
$ws.Range("H129").Value = 1144.9445
$ws.Range("I129").Value = 564.44446
$ws.Range("J129").Value = 1725.4445
$ws.Range("K129").Value = 1693.33338
$ws.Range("L129").Value = 5176.333500000001
$ws.Range("M129").Value = 3306.66662
$ws.Range("N129").Value = -15176.3335
$ws.Range("H130").Value = 2617.5
$ws.Range("I130").Value = 2434.2856
$ws.Range("J130").Value = 3900
$ws.Range("K130").Value = 7302.8568
$ws.Range("L130").Value = 11700
$ws.Range("M130").Value = -2282.8568
$ws.Range("N130").Value = -21740
$ws.Range("H131").Value = 2573.6064
$ws.Range("J131").Value = 1766.2264
$ws.Range("L131").Value = 5298.6792
$ws.Range("N131").Value = -15378.6792
$ws.Range("H132").Value = 166668260
$ws.Range("I132").Value = 333333660
$ws.Range("J132").Value = 2833.3333
$ws.Range("K132").Value = 3000002940
$ws.Range("L132").Value = 25499.9997
$ws.Range("M132").Value = -3000000410
$ws.Range("N132").Value = -30559.9997
$ws.Range("H133").Value = 8033.3335
$ws.Range("I133").Value = 4100
$ws.Range("K133").Value = 12300
$ws.Range("M133").Value = -7240
$ws.Range("H134").Value = 3085.9092
$ws.Range("I134").Value = 2783.6843
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 8351.052899999999
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -3281.052899999999
$ws.Range("N134").Value = -25140
$ws.Range("H136").Value = 8432.666999999999
$ws.Range("I136").Value = 530
$ws.Range("J136").Value = 10013.2
$ws.Range("K136").Value = 1590
$ws.Range("L136").Value = 30039.6
$ws.Range("M136").Value = 3510
$ws.Range("N136").Value = -40239.60000000001
$ws.Range("H137").Value = 18324.4
$ws.Range("I137").Value = 2825
$ws.Range("J137").Value = 22199.25
$ws.Range("K137").Value = 8475
$ws.Range("L137").Value = 66597.75
$ws.Range("M137").Value = -3375
$ws.Range("N137").Value = -76797.75
$ws.Range("H138").Value = 1870.2307
$ws.Range("I138").Value = 1773.3334
$ws.Range("J138").Value = 3033
$ws.Range("K138").Value = 5320.0002
$ws.Range("L138").Value = 9099
$ws.Range("M138").Value = -180.0002000000004
$ws.Range("N138").Value = -19379
$ws.Range("H139").Value = 2503.56
$ws.Range("I139").Value = 1227.1666
$ws.Range("J139").Value = 5785.7144
$ws.Range("K139").Value = 3681.4998
$ws.Range("L139").Value = 17357.1432
$ws.Range("M139").Value = 1458.5002
$ws.Range("N139").Value = -27637.1432
$ws.Range("H140").Value = 1690.5
$ws.Range("I140").Value = 904.3889
$ws.Range("J140").Value = 4048.8333
$ws.Range("K140").Value = 2713.1667
$ws.Range("L140").Value = 12146.4999
$ws.Range("M140").Value = 2466.8333
$ws.Range("N140").Value = -22506.4999
$ws.Range("H141").Value = 136626.8
$ws.Range("I141").Value = 145790.64
$ws.Range("J141").Value = 8333
$ws.Range("K141").Value = 437371.92
$ws.Range("L141").Value = 24999
$ws.Range("M141").Value = -432191.92
$ws.Range("N141").Value = -35359

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 17545484
$ws.Range("I132").Value = 23257258
$ws.Range("J132").Value = 2176.2144
$ws.Range("K132").Value = 69771774
$ws.Range("L132").Value = 6528.6432
$ws.Range("M132").Value = -69769244
$ws.Range("N132").Value = -11588.6432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1900.8
$ws.Range("I61").Value = 1802
$ws.Range("J61").Value = 1966.6666
$ws.Range("K61").Value = 1802
$ws.Range("L61").Value = 1966.6666
$ws.Range("M61").Value = -1600
$ws.Range("N61").Value = -2370.6666
$ws.Range("H100").Value = 3263.4707
$ws.Range("I100").Value = 1728
$ws.Range("J100").Value = 4990.875
$ws.Range("K100").Value = 1728
$ws.Range("L100").Value = 4990.875
$ws.Range("M100").Value = -1187
$ws.Range("N100").Value = -6072.875
$ws.Range("H113").Value = 1900.8
$ws.Range("I113").Value = 1802
$ws.Range("J113").Value = 1966.6666
$ws.Range("K113").Value = 1802
$ws.Range("L113").Value = 1966.6666
$ws.Range("M113").Value = 368
$ws.Range("N113").Value = -6306.6666
$ws.Range("H132").Value = 3348.283
$ws.Range("I132").Value = 3249.62
$ws.Range("K132").Value = 9748.860000000001
$ws.Range("M132").Value = -7218.860000000001
$ws.Range("H136").Value = 922.2895
$ws.Range("I136").Value = 556.2258
$ws.Range("J136").Value = 2543.4285
$ws.Range("K136").Value = 1668.6774
$ws.Range("L136").Value = 7630.2855
$ws.Range("M136").Value = 881.3226
$ws.Range("N136").Value = -12730.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7133759
$ws.Range("I132").Value = 7846774.5
$ws.Range("J132").Value = 3599.6
$ws.Range("K132").Value = 23540323.5
$ws.Range("L132").Value = 10798.8
$ws.Range("M132").Value = -23537793.5
$ws.Range("N132").Value = -15858.8
$ws.Range("H136").Value = 2920618
$ws.Range("I136").Value = 7485.3
$ws.Range("K136").Value = 22455.9
$ws.Range("M136").Value = -19905.9
